# Auto-generated edit script: apply cryptos.xlsx price/volume refresh (commit: "Updated symbol list on Tue Jan 10 09:32:50 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: set a cell to a literal text value, preventing Excel from
# auto-converting numeric-/percent-looking strings into numbers.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# --- Price (D) / Volume(1h) (E) updates ---
Set-TextCell "D2" "275.38"
Set-TextCell "E2" "-1.11%"
Set-TextCell "D3" "26.51"
Set-TextCell "E3" "-2.83%"
Set-TextCell "D4" "4.892"
Set-TextCell "E4" "2.10%"
Set-TextCell "D5" "0.06341"
Set-TextCell "E5" "0.85%"
Set-TextCell "D6" "6.893"
Set-TextCell "E6" "-0.44%"
Set-TextCell "D7" "3.312"
Set-TextCell "E7" "1.28%"
Set-TextCell "D8" "1.290"
Set-TextCell "E8" "36.78%"
Set-TextCell "D9" "0.8670"
Set-TextCell "E9" "-1.26%"
Set-TextCell "D10" "0.1537"
Set-TextCell "E10" "5.35%"
Set-TextCell "D11" "0.05029"
Set-TextCell "E11" "-2.53%"
Set-TextCell "D12" "0.07400"
Set-TextCell "E12" "1.70%"
Set-TextCell "D13" "0.02939"
Set-TextCell "E13" "-5.17%"
Set-TextCell "D14" "0.09047"
Set-TextCell "E14" "-0.24%"
Set-TextCell "D15" "0.001574"
Set-TextCell "E15" "1.12%"
Set-TextCell "D16" "0.0006319"
Set-TextCell "E16" "0.50%"
Set-TextCell "D17" "0.005911"
Set-TextCell "E17" "0.71%"
Set-TextCell "D18" "3.448"
Set-TextCell "E18" "-0.08%"
Set-TextCell "D19" "2.272"
Set-TextCell "E19" "-0.55%"
Set-TextCell "E21" "0.94%"
Set-TextCell "D22" "3.894"
Set-TextCell "E22" "1.16%"
Set-TextCell "D23" "0.04369"
Set-TextCell "E23" "0.78%"
Set-TextCell "D24" "0.001169"
Set-TextCell "E24" "-1.00%"
Set-TextCell "D25" "0.004254"
Set-TextCell "E25" "-0.70%"
Set-TextCell "E26" "-0.14%"
Set-TextCell "D27" "0.0001677"
Set-TextCell "E27" "-0.81%"
Set-TextCell "D40" "0.04109"
Set-TextCell "E40" "1.04%"
Set-TextCell "D41" "0.006984"
Set-TextCell "E41" "6.02%"
Set-TextCell "E42" "1.18%"
Set-TextCell "E43" "-1.42%"
Set-TextCell "D44" "0.01080"
Set-TextCell "E44" "-8.84%"
Set-TextCell "D45" "0.00005266"
Set-TextCell "E45" "2.54%"

# --- Rows 46/47 swapped (BOLO <-> CoinbaseStockToken) with refreshed values ---
$ws.Range("B46").Value = "CoinbaseStockToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell "D46" "0.01998"
Set-TextCell "E46" "-11.24%"

$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell "D47" "1.490"
Set-TextCell "E47" "-37.35%"
